$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values - B2 and D2 updated, C2 and E2 cleared
$ws.Range("B2").Value = 13.529346766142577
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 11.262458768517522
$ws.Range("E2").ClearContents()

# Row 3 values (B3:E3)
$ws.Range("B3").Value = 8.4170495821710336
$ws.Range("C3").Value = -16.827364875250563
$ws.Range("D3").Value = 5.9900348240649164
$ws.Range("E3").Value = -19.35657950153503

# Update selection to match new active range
$ws.Range("B1:E3").Select()
